$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (B7="M") moves up to become row 6 (B6="M"); clear the old row 7
$ws.Range("B7").ClearContents()
$ws.Range("B6").Value = "M"

# Old row 8 (A8="Div", B8="M", C8=632589) becomes new row 9, unchanged values
$ws.Range("A9").Value = "Div"
$ws.Range("B9").Value = "M"
$ws.Range("C9").Value = 632589

# Old A8 is cleared (no longer present); old row 8 now only keeps B8="M" and C8=63258
$ws.Range("A8").ClearContents()
$ws.Range("B8").Value = "M"
$ws.Range("C8").Value = 63258

# Update the active selection to match the new state
$ws.Range("E7").Select()
